# Add data for 2022-07-31
# - Rename the sheet / update title text from "...07-22" to "...07-23"
# - Update July (row 8) 2022-through-date value: 125 -> 133
# - Update Total (row 14) 2022-through-date value: 931 -> 939

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet tab name
$ws.Name = "Through 2022-07-23"

# Header label in column I (row 1) - shared string "2022 (through 07-22)" -> "2022 (through 07-23)"
$ws.Range("I1").Value = "2022 (through 07-23)"

# July row, 2022 column
$ws.Range("I8").Value = 133

# Total row, 2022 column
$ws.Range("I14").Value = 939
